{"js": "// Apply the Technical Report text edits described by the diff.\n// (The diff's many <w:proofErr> insertions are Word's own live spell/grammar\n// -checker artifacts from re-typing/re-saving the paragraphs; they carry no\n// textual content. The actual substantive change is the wording below.)\n\nasync function replaceOnce(searchText, newText, matchCase) {\n  const results = context.document.body.search(searchText, {\n    matchCase: matchCase !== false,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"...original csv and json formatting into data frames...\"\n//    -> \"...original .csv and .json formatting into dataframes...\"\nawait replaceOnce(\n  \"original csv and json formatting into data frames\",\n  \"original .csv and .json formatting into dataframes\"\n);\n\n// 2. \"Pytho\" + \"n and SQLalchemy:\"  -> \"...SQLAlchemy:\" (capitalize the A)\nawait replaceOnce(\"SQLalchemy\", \"SQLAlchemy\");\n\n// 3. First \"data frames\" in the to_csv/to_sql paragraph becomes \"dataframes\"\nawait replaceOnce(\n  \"Using the data frames we created\",\n  \"Using the dataframes we created\"\n);\n\n// 4. Third \"data frames\" in the same paragraph becomes \"dataframes\"\n//    (the middle occurrence, \"export and read the data frames\", is left as-is)\nawait replaceOnce(\n  \"push the data frames and create\",\n  \"push the dataframes and create\"\n);\n\n// 5. \"...latitude and longitude of a U.S. fireball sighting.\"\n//    -> \"...latitude and longitude of the sighting.\"\nawait replaceOnce(\n  \"of a U.S. fireball sighting\",\n  \"of the sighting\"\n);\n", "ps1": "# Apply the Technical Report text edits described by the diff.\n# (The diff's many <w:proofErr> insertions are Word's own live spell/grammar\n# -checker artifacts from re-typing/re-saving the paragraphs; they carry no\n# textual content. The actual substantive change is the wording below.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n    if (-not $found) {\n        throw \"Text not found: $find\"\n    }\n}\n\n# 1. \"...original csv and json formatting into data frames...\"\n#    -> \"...original .csv and .json formatting into dataframes...\"\nReplace-Text \"original csv and json formatting into data frames\" \"original .csv and .json formatting into dataframes\"\n\n# 2. \"Pytho\" + \"n and SQLalchemy:\"  -> \"...SQLAlchemy:\" (capitalize the A)\nReplace-Text \"SQLalchemy\" \"SQLAlchemy\"\n\n# 3. First \"data frames\" in the to_csv/to_sql paragraph becomes \"dataframes\"\nReplace-Text \"Using the data frames we created\" \"Using the dataframes we created\"\n\n# 4. Third \"data frames\" in the same paragraph becomes \"dataframes\"\n#    (the middle occurrence, \"export and read the data frames\", is left as-is)\nReplace-Text \"push the data frames and create\" \"push the dataframes and create\"\n\n# 5. \"...latitude and longitude of a U.S. fireball sighting.\"\n#    -> \"...latitude and longitude of the sighting.\"\nReplace-Text \"of a U.S. fireball sighting\" \"of the sighting\"\n"}
